# Weekly update: a new sampling row is inserted at row 8 (pushing the
# existing historical rows down by one), representing the newest
# "Fruta / hortaliza, semanal" observation for Coco at Vega Central
# Mapocho de Santiago.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 8; rows 8:42 shift down to 9:43,
# inheriting the formatting (incl. the date style on column D) of the
# row that used to be at 8.
$ws.Rows.Item(8).Insert()

# Populate the freshly inserted row 8 with the new weekly observation.
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 45250
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = "Fruta"
$ws.Range("G8").Value = 100108
$ws.Range("H8").Value = "Tropicales y subtropicales"
$ws.Range("I8").Value = 100108007
$ws.Range("J8").Value = "Coco"
$ws.Range("K8").Value = "Sin especificar"
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 28000
$ws.Range("O8").Value = 28000
$ws.Range("P8").Value = 28000
$ws.Range("Q8").Value = "$/malla 20 unidades"
$ws.Range("R8").Value = "Perú"
$ws.Range("S8").Value = 1400
$ws.Range("T8").Value = 20
